# Completion of Give Feedback, Checking Firefox execution issues when
# running the master test suite 200819
#
# The "Plumber in Rugby" test row is replaced with a "Blocked Sinks in
# Rugby" scenario (tradeLocationVerification column G, row 5), and the
# active selection is moved back to cell E1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = "Blocked Sinks in Rugby"

$null = $ws.Range("E1").Select()
